$d = $word.ActiveDocument

# Locate the paragraph that ends with "Domain Models Entities / Relationships...";
# new paragraphs are inserted right after it (before the existing blank
# paragraph that precedes "Backend Architecture (To Do):").
$rng = $d.Content
$found = $rng.Find.Execute(
    "Domain Models Entities / Relationships: transforms of underlying Entities given Relationships contents.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor paragraph not found"
}

$rng.Collapse(0)
$anchorIndex = $d.Range($rng.Start, $rng.Start).Paragraphs(1).Index

# New paragraphs to add, in order; "" means an empty (blank) paragraph.
$newParas = @(
    "",
    "Relationships (upper domain): before, during, after, cause, effect, implies, partOf, etc.",
    "",
    "Relationsip assertions reified / parsed as / from Relation Statements (Messages predicates).",
    "",
    "Entity Relationships rendered / parsed as / from Relations, Kinds, Statement, Resources (Message contents)."
)

$currentIndex = $anchorIndex
foreach ($text in $newParas) {
    $p = $d.Paragraphs($currentIndex)
    $p.Range.InsertParagraphAfter()
    $currentIndex = $currentIndex + 1
    if ($text -ne "") {
        $newP = $d.Paragraphs($currentIndex)
        $newP.Range.Text = $text
    }
}
